$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-21 Tuesday" "2025-10-22 Wednesday"

Replace-Text "369÷7=52, 5" "649÷5=129, 4"
Replace-Text "633÷5=126, 3" "278÷8=34, 6"
Replace-Text "406÷6=67, 4" "565÷8=70, 5"
Replace-Text "423÷6=70, 3" "755÷2=377, 1"
Replace-Text "177÷8=22, 1" "663÷5=132, 3"

Replace-Text "674÷2=337, 0" "443÷6=73, 5"
Replace-Text "553÷2=276, 1" "897÷4=224, 1"
Replace-Text "373÷7=53, 2" "939÷7=134, 1"
Replace-Text "312÷6=52, 0" "296÷4=74, 0"
Replace-Text "903÷3=301, 0" "257÷9=28, 5"

Replace-Text "404÷7=57, 5" "696÷2=348, 0"
Replace-Text "641÷9=71, 2" "491÷7=70, 1"
Replace-Text "433÷5=86, 3" "994÷6=165, 4"
Replace-Text "786÷9=87, 3" "386÷3=128, 2"
Replace-Text "588÷7=84, 0" "841÷3=280, 1"

Replace-Text "519÷9=57, 6" "408÷8=51, 0"
Replace-Text "157÷2=78, 1" "228÷2=114, 0"
Replace-Text "320÷6=53, 2" "679÷6=113, 1"
Replace-Text "706÷8=88, 2" "252÷6=42, 0"
Replace-Text "716÷6=119, 2" "622÷7=88, 6"

Replace-Text "136÷6=22, 4" "432÷9=48, 0"
Replace-Text "577÷4=144, 1" "541÷5=108, 1"
Replace-Text "455÷4=113, 3" "836÷7=119, 3"
Replace-Text "302÷2=151, 0" "469÷9=52, 1"
Replace-Text "948÷8=118, 4" "571÷8=71, 3"

$d.Save()
